$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.244.18'
$ws.Range('E2').Value = '  +0.05%  '

$ws.Range('D3').Value = '3.982.42'
$ws.Range('E3').Value = '  -1.65%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = "'609.35"
$ws.Range('E5').Value = '  +6.02%  '

$ws.Range('D6').Value = "'172.46"
$ws.Range('E6').Value = '  +13.49%  '

$ws.Range('D7').Value = "'0.703"
$ws.Range('E7').Value = '  +1.15%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').Value = "'0.802"
$ws.Range('E9').Value = '  +4.86%  '

$ws.Range('E10').Value = '  +8.59%  '

$ws.Range('D11').Value = "'57.03"
$ws.Range('E11').Value = '  +5.74%  '

$ws.Range('D12').Value = "'0.0000337"
$ws.Range('E12').Value = '  +3.24%  '

$ws.Range('D13').Value = "'11.82"
$ws.Range('E13').Value = '  +5.70%  '

$ws.Range('D14').Value = '4.615.10'
$ws.Range('E14').Value = '  -1.73%  '

$ws.Range('D15').Value = '3.987.84'
$ws.Range('E15').Value = '  -1.71%  '

$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').Value = "'14.36"
$ws.Range('E16').Value = '  +0.04%  '

$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').Value = "'1.26"
$ws.Range('E17').Value = '  +2.43%  '

$ws.Range('D18').Value = "'21.13"
$ws.Range('E18').Value = '  +0.99%  '

$ws.Range('D19').Value = '73.111.56'
$ws.Range('E19').Value = '  -0.07%  '

$ws.Range('E20').Value = '  -0.97%  '

$ws.Range('D21').Value = "'469.19"
$ws.Range('E21').Value = '  +5.55%  '

$ws.Range('D22').Value = "'4.83"
$ws.Range('E22').Value = '  +7.20%  '

$ws.Range('D23').Value = "'97.73"
$ws.Range('E23').Value = '  -0.65%  '

$ws.Range('E24').Value = '  -4.20%  '

$ws.Range('D25').Value = "'14.37"
$ws.Range('E25').Value = '  -2.09%  '

$ws.Range('D26').Value = "'4.25"
$ws.Range('E26').Value = '  -0.42%  '

$ws.Range('D27').Value = "'11.32"
$ws.Range('E27').Value = '  -0.57%  '

$ws.Range('D28').Value = "'10.71"
$ws.Range('E28').Value = '  -3.18%  '

$ws.Range('E29').Value = '  -1.04%  '

$ws.Range('D30').Value = "'36.64"
$ws.Range('E30').Value = '  -1.36%  '

$ws.Range('D31').Value = "'8.07"
$ws.Range('E31').Value = '  +2.52%  '

$ws.Range('D32').Value = "'14.16"
$ws.Range('E32').Value = '  +3.36%  '

$ws.Range('D33').Value = "'50.25"
$ws.Range('E33').Value = '  +3.32%  '

$ws.Range('D34').Value = "'0.130"
$ws.Range('E34').Value = '  -2.77%  '

$ws.Range('E35').Value = '  +14.04%  '

$ws.Range('D36').Value = "'70.69"
$ws.Range('E36').Value = '  +2.83%  '

$ws.Range('D37').Value = "'642.96"
$ws.Range('E37').Value = '  -7.64%  '

$ws.Range('D38').Value = "'0.437"
$ws.Range('E38').Value = '  -2.70%  '

$ws.Range('E39').Value = '  -0.93%  '

$ws.Range('E40').Value = '  +0.83%  '

$ws.Range('E41').Value = '  -0.04%  '

$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = "'3.30"
$ws.Range('E42').Value = '  +43.66%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  -0.04%  '

$ws.Range('D44').Value = "'0.0487"
$ws.Range('E44').Value = '  -1.68%  '

$ws.Range('D45').Value = "'10.59"
$ws.Range('E45').Value = '  -6.88%  '

$ws.Range('D46').Value = "'0.152"
$ws.Range('E46').Value = '  -0.01%  '

$ws.Range('D47').Value = "'3.01"
$ws.Range('E47').Value = '  -9.18%  '

$ws.Range('D48').Value = "'0.000303"
$ws.Range('E48').Value = '  +10.26%  '

$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').Value = "'2.66"
$ws.Range('E49').Value = '  -3.34%  '

$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = "'3.45"
$ws.Range('E50').Value = '  +3.94%  '

$ws.Range('D51').Value = '2.827.00'
$ws.Range('E51').Value = '  +1.23%  '

